$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 51, pushing the existing rows 51-54 down to 53-56.
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(51).Insert()

# New row 51: Agrícola del Norte S.A. de Arica - Ají - Cristal - Primera
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 44578
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = 100112021
$ws.Range("G51").Value = "Ají"
$ws.Range("H51").Value = "Cristal"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 120
$ws.Range("K51").Value = 16000
$ws.Range("L51").Value = 17000
$ws.Range("M51").Value = 16500
$ws.Range("N51").Value = "$/caja 15 kilos"
$ws.Range("O51").Value = "Región de Arica y Parinacota"
$ws.Range("P51").Value = 1100
$ws.Range("Q51").Value = 15
$ws.Range("R51").Value = "Hortaliza"

# New row 52: Agrícola del Norte S.A. de Arica - Ají - Inferno - Primera
$ws.Range("A52").Value = 1
$ws.Range("B52").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C52").Value = "Arica y Parinacota"
$ws.Range("D52").Value = 44578
$ws.Range("E52").Value = 15
$ws.Range("F52").Value = 100112021
$ws.Range("G52").Value = "Ají"
$ws.Range("H52").Value = "Inferno"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 130
$ws.Range("K52").Value = 18000
$ws.Range("L52").Value = 20000
$ws.Range("M52").Value = 19000
$ws.Range("N52").Value = "$/caja 15 kilos"
$ws.Range("O52").Value = "Región de Arica y Parinacota"
$ws.Range("P52").Value = 1267
$ws.Range("Q52").Value = 15
$ws.Range("R52").Value = "Hortaliza"
